$d = $word.ActiveDocument

# 1) mobilizedResource -> resourceInfo
$d.Content.Find.Execute("mobilizedResource", $true, $false, $false, $false, $false,
                         $true, 1, $false, "resourceInfo", 2) | Out-Null

# 2) "Ressource engagée / à engager" -> "Ressource"
$d.Content.Find.Execute("Ressource engagée / à engager", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Ressource", 2) | Out-Null

# 3) Merge the two-sentence description (separated by a manual line break)
#    into a single sentence, dropping the " pour le message RS-RI" suffix
#    and the entire second sentence about RS-ER.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Text = "Objet permettant de communquer la liste des ressource et vecteurs mobilisés en 15-15 et 15-SMUR pour le message RS-RI^lObjet permettant de communiquer la liste des ressources à engager en 15-SMUR pour le message RS-ER"
$find.Replacement.Text = "Objet permettant de communquer la liste des ressource et vecteurs mobilisés en 15-15 et 15-SMUR"
$find.Forward = $true
$find.Wrap = 1
$find.Format = $false
$find.MatchCase = $false
$find.MatchWholeWord = $false
$find.MatchWildcards = $false
$find.MatchSoundsLike = $false
$find.MatchAllWordForms = $false
$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null
